# Generate Report for Handoff
# This script updates the localization-status workbook to reflect the
# "Ready for handoff" report generation: status text flips from
# "Handed back: in sync with en-US" to "Ready for handoff", timestamps are
# refreshed, the zh-cn priority changes from "ht" to "mt", and an error
# detail message is recorded for the 6e24cf68 file (stale handback version).

$wb = $excel.ActiveWorkbook

$statusOld = "Handed back: in sync with en-US"
$statusNew = "Ready for handoff"

$overviewDateOld = "2016-11-29 05:29:55"
$overviewDateNew = "2016-11-29 05:31:59"

$zhHandoffDateOld = "2016-11-29 05:29:36"
$zhHandoffDateNew = "2016-11-29 05:31:45"

$priorityOld = "ht"
$priorityNew = "mt"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/798141b252550f36f5813e352255cc5acf5f813b/e2e/6e24cf68-6a8c-485c-afad-9fba52202759.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/29b42280ea3730cf874c8a4f0b8cb1a1c714d654/e2e/6e24cf68-6a8c-485c-afad-9fba52202759.md."

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("G2").Value = $overviewDateNew

$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Range("G3").Value = $overviewDateNew

$wsOverview.Columns.Item(5).ColumnWidth = 16.38265482584637
$wsOverview.Columns.Item(6).ColumnWidth = 16.38265482584637

# ----- zh-cn sheet -----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("E2").Value = $priorityNew
$wsZh.Range("H2").Value = $zhHandoffDateNew
$wsZh.Range("P2").Value = $errorDetail

$wsZh.Range("C3").Value = $statusNew
$wsZh.Range("E3").Value = $priorityNew
$wsZh.Range("H3").Value = $zhHandoffDateNew

$wsZh.Columns.Item(3).ColumnWidth = 16.38265482584637
$wsZh.Columns.Item(16).ColumnWidth = 39.166666666666664

# ----- de-de sheet -----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("E2").Value = $priorityNew
# H2/H3 on de-de shared the same underlying string as the Overview's
# "Latest HO Xliff Generate Date" (both were "2016-11-29 05:29:55"), so it
# moves in lockstep with the Overview date refresh.
$wsDe.Range("H2").Value = $overviewDateNew
$wsDe.Range("P2").Value = $errorDetail

$wsDe.Range("C3").Value = $statusNew
$wsDe.Range("E3").Value = $priorityNew
$wsDe.Range("H3").Value = $overviewDateNew

$wsDe.Columns.Item(3).ColumnWidth = 16.38265482584637
$wsDe.Columns.Item(16).ColumnWidth = 39.166666666666664
